$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C32").Value = 50.55
$ws.Range("E33").Value = 66.5
$ws.Range("D48").Value = 1.2
$ws.Range("C59").Value = 13.35
$ws.Range("C94").Value = 13.35
$ws.Range("D110").Value = 1.2
$ws.Range("C114").Value = 13.35
$ws.Range("D125").Value = 1.2
$ws.Range("C128").Value = 13.35
$ws.Range("D129").Value = 1.2
$ws.Range("C130").Value = 50.55
$ws.Range("D136").Value = 1.2
$ws.Range("D148").Value = 1.2
$ws.Range("D152").Value = 1.2
